$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that used to sit right after the
#    "...Tests" run (Word leaves one of these behind at the last edit
#    point; it is being cleaned up here).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    [void]$d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Locate the "Запоминаем последний номер страницы в состоянии"
#    heading paragraph (Lesson 5-4) so three new paragraphs can be
#    inserted right before it, introducing a new "Дозаполнение
#    сведениями" page/title section.
# ---------------------------------------------------------------------
function Find-Target {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith("Запоминаем последний номер страницы в состоянии")) {
            return $p
        }
    }
    return $null
}

$target = Find-Target

if ($target -ne $null) {
    $ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

    # Create three blank paragraphs right before the target, re-locating
    # the target paragraph after each insertion so stale paragraph
    # references/offsets from before the edit are never relied upon.
    # Because the target paragraph carries the "heading 3" / numId 7
    # (bulleted) paragraph formatting, each freshly split-off paragraph
    # inherits that same formatting to start with.
    [void]$target.Range.InsertParagraphBefore()
    $target = Find-Target
    [void]$target.Range.InsertParagraphBefore()
    $target = Find-Target
    [void]$target.Range.InsertParagraphBefore()
    $target = Find-Target

    $p3 = $target.Previous()    # empty paragraph holding the "_GoBack" bookmark
    $p2 = $p3.Previous()        # new body text
    $p1 = $p2.Previous()        # new heading: "Дозаполнение сведениями"

    # Paragraph 1 keeps the inherited "heading 3" + list formatting.
    $xml1 = "<w:p $ns>" +
            "<w:pPr><w:pStyle w:val=`"3`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr></w:pPr>" +
            "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Дозаполнение</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
            "<w:r><w:t xml:space=`"preserve`"> сведениями</w:t></w:r>" +
            "</w:p>"
    [void]$p1.Range.InsertXML($xml1)

    # Paragraph 2 is plain body text (no heading/list formatting).
    $xml2 = "<w:p $ns>" +
            "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Дозаполняем</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
            "<w:r><w:t xml:space=`"preserve`"> модельку с квартирой адресом дома, макс. Этажом. Дома заполняем кол-вом квартир.</w:t></w:r>" +
            "</w:p>"
    [void]$p2.Range.InsertXML($xml2)

    # Paragraph 3 is empty apart from the relocated "_GoBack" bookmark.
    $xml3 = "<w:p $ns><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
    [void]$p3.Range.InsertXML($xml3)
}
